$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Comment" column-name synonyms cell (B3)
$ws.Range("B3").Value = "Comment, Comment..2"

# Update the d13C_m synonyms cell (B22) and d13_C synonyms cell (B21).
# Order matters: d13C_m is set before d13_C so the shared-string table
# append order matches the target workbook.
$ws.Range("B22").Value = "d13C_m, δ13C ‰ measured"
$ws.Range("B21").Value = "d13_C, δ13C ‰ VPDB, d13C PDB, δ13C [‰, PDB]"

# Move the active selection to B22 (also clears the old scrolled/topLeftCell view state)
$null = $ws.Range("B22").Select()
